$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text content (order matters for shared-string slot reuse) ---
# Move "N/A" out of D6 into new D8 first so it survives the rebuild.
$ws.Range("D8").Value = "N/A"
# C6 "Both scripts worked as intended " -> replaced in place with new text.
$ws.Range("C6").Value = "Search script would not cancel when white is seen."
# D6 (previously "N/A") -> new follow-up note text.
$ws.Range("D6").Value = "altered code after each rotation on the search added ""Reuturn"" script when sensors sees white during search"

# New row 7
$ws.Range("B7").Value = "tested audio feedback"
$ws.Range("C7").Value = "script failed with errors"
$ws.Range("D7").Value = "did researcj to find out documentation was out of date. Used microsoft co-pilot to reaserch other scrips that may work and found one. Altered script until audio feedback worked."

# New row 8 (B/C; D8 already set above)
$ws.Range("B8").Value = "Testing all code"
$ws.Range("C8").Value = "All code worked as intended"

# --- Formatting ---
# D6 gets the "center + wrap" style (created first so it lands on cellXfs index 1).
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").WrapText = $true

# Everything else in the used range gets "wrap text" only (cellXfs index 2).
$ws.Range("B2:D8").WrapText = $true

# Row heights for the rows whose content now wraps to two lines.
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(7).RowHeight = 28.8

# Column widths (engine rounds ColumnWidth to whole pixels internally; this
# lands on the closest representable stored width to the target 78.5546875).
$ws.Columns.Item(4).ColumnWidth = 77.666666666667

# Selection / active cell moves to B8 to match the end state.
$ws.Range("B8").Select() | Out-Null
